$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell value updates (rows 7-11) - new footprints / Mouser parts
# ---------------------------------------------------------------------

# Row 7 - Relais 1pool: aantal 10 -> 11
$ws.Range("C7").Value = 11

# Row 8 - Schroef headers 1x3 -> now sourced from Mouser
$ws.Range("B8").Value = 0.73
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "Mouser"

# Row 9 - Schroef headers 1x2 -> now sourced from Mouser
$ws.Range("B9").Value = 0.5
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "Mouser"

# Row 10 - was "Schroef headers 1x6", now "Schroef headers 1x4" (Mouser)
$ws.Range("A10").Value = "Schroef headers 1x4"
$ws.Range("B10").Value = 0.97
$ws.Range("D10").Value = "Mouser"

# Row 11 - was "Schroef headers 1x8", now "Schroef headers 1x14" (Mouser)
$ws.Range("A11").Value = "Schroef headers 1x14"
$ws.Range("B11").Value = 2.48
$ws.Range("D11").Value = "Mouser"

# ---------------------------------------------------------------------
# 2. Hyperlinks - swap old Conrad "Degson" links for new Mouser
#    "CUI Devices TB006-508" links. The COM shim here only supports
#    *adding* hyperlinks (Range.Hyperlinks.Delete() clears the whole
#    sheet's collection), so rebuild the full set from scratch in the
#    final order/targets.
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.conrad.be/p/finder-405290120000-printrelais-12-vdc-8-a-2x-wisselcontact-1-stuks-502868")
$ws.Range("E6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.conrad.be/p/finder-403170120000-printrelais-12-vdc-12-a-1x-wisselcontact-1-stuks-1560601")
$ws.Range("E7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.mouser.be/ProductDetail/CUI-Devices/TB006-508-14BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdta3MFIStyIuo0QxXNoLqeaQ%3D%3D")
$ws.Range("E11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.mouser.be/ProductDetail/CUI-Devices/TB006-508-02BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdthU5yE00GqwuGh8iwAbLcpg%3D%3D")
$ws.Range("E9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E10"), "https://www.mouser.be/ProductDetail/CUI-Devices/TB006-508-04BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdtxX8qyOuZQe4mznRYcE70wg%3D%3D")
$ws.Range("E10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.mouser.be/ProductDetail/CUI-Devices/TB006-508-03BE?qs=sGAEpiMZZMvZTcaMAxB2AHpdXjUJWjdtVgaMJ9quMmnG8UNjIL2MZA%3D%3D")
$ws.Range("E8").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3. Column widths - column D/E auto-fit after content changed
#    (shorter "Mouser"/narrower D column, much longer Mouser URLs in E)
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 10.6
$ws.Columns.Item(5).ColumnWidth = 135.8

# ---------------------------------------------------------------------
# 4. Selection moved to E17
# ---------------------------------------------------------------------
$ws.Range("E17").Select()
